# This workbook's "Fecha" (date, col D) together with the quality/volume/price/
# origin/unit columns (L..T) for each data row (rows 2-25) were re-shuffled
# between rows (the descriptive columns A,B,C,E..K stay the same for every
# row). Build the new row -> source row mapping and copy the values over.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# newRow -> sourceRow (which row's D,L,M,N,O,P,Q,R,S,T values move into newRow)
$rowMap = [ordered]@{
    2  = 8
    3  = 4
    4  = 16
    5  = 24
    6  = 12
    7  = 10
    8  = 19
    9  = 2
    10 = 15
    11 = 9
    12 = 17
    13 = 21
    14 = 20
    15 = 25
    16 = 3
    17 = 11
    18 = 14
    19 = 5
    20 = 6
    21 = 13
    22 = 7
    23 = 18
    24 = 22
    25 = 23
}

# Columns whose contents move together as a group, keyed by column letter.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot every source row's values BEFORE writing anything, since this is a
# permutation (values are swapped among rows, not just copied one-way).
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the shuffled values back using the snapshot as the source of truth.
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $srcVals[$col]
    }
}
